$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 45) mirroring the existing table's layout.
$newRow = 45

$ws.Cells.Item($newRow, 1).Value = 46000
$ws.Cells.Item($newRow, 2).Value = 5620
$ws.Cells.Item($newRow, 3).Value = 4321
$ws.Cells.Item($newRow, 4).Value = 3988
$ws.Cells.Item($newRow, 5).Value = 239
$ws.Cells.Item($newRow, 6).Value = 60
$ws.Cells.Item($newRow, 7).Value = 30
$ws.Cells.Item($newRow, 8).Value = 3
$ws.Cells.Item($newRow, 9).Value = 1

# Match the date number formatting used by the column above (style s="4").
$ws.Range("A45").NumberFormat = $ws.Range("A44").NumberFormat

# Move/extend the active selection to the newly added row, like Excel does
# right after data entry.
$ws.Range("A45:I45").Select()
